$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.418.48"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.655.07"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.535"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.34%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("E11").Value = "  +3.50%  "

$ws.Range("D12").Value = "1.889.03"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").Value = "1.656.61"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "27.413.53"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.34%  "

$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.69%  "

$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("E26").Value = "  -1.36%  "

$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("E31").Value = "  -4.25%  "

$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.420.54"
$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  -2.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.567"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.76%  "

$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("D46").Value = "1.798.18"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("E49").Value = "  -3.70%  "

$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
